{"js": "// Add a new \"Utah\" paragraph at the very start of the document body, then\n// prepend the Saturday-weather sentence onto what was originally the first\n// (and only) paragraph \u2014 which holds the _GoBack bookmark and must keep it.\n\nconst body = context.document.body;\n\n// 1) Insert a brand-new paragraph \"Utah\" before everything else.\nbody.insertParagraph(\"Utah\", \"Start\");\n\n// 2) Re-fetch the paragraph list; the original (bookmarked) paragraph is now\n//    the second paragraph, and we insert the weather sentence at its start.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst originalParagraph = paragraphs.items[1];\noriginalParagraph.insertText(\n  \"On Saturday we woke up to several inches of heavy wet snow on the ground. \",\n  \"Start\"\n);\n\nawait context.sync();\n", "ps1": "# Add a new \"Utah\" paragraph at the very start of the document, then\n# prepend the Saturday-weather sentence onto what was originally the first\n# (and only) paragraph -- which carries the _GoBack bookmark and must keep it.\n\n$d = $word.ActiveDocument\n\n# 1) Insert a brand-new empty paragraph before the original first paragraph,\n#    then give it the text \"Utah\".\n$firstParagraph = $d.Paragraphs(1)\n$firstParagraph.Range.InsertParagraphBefore()\n$d.Paragraphs(1).Range.Text = \"Utah\"\n\n# 2) The original (bookmarked) paragraph is now the second paragraph; insert\n#    the weather sentence at its start, ahead of the bookmark.\n$d.Paragraphs(2).Range.InsertBefore(\"On Saturday we woke up to several inches of heavy wet snow on the ground. \")\n"}
